$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing GPIO pin numbers (column C, rows 2-15) ---
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 9
$ws.Range("C4").Value = 13
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 11
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 2
$ws.Range("C11").Value = 3
$ws.Range("C12").Value = 4
$ws.Range("C13").Value = 5
$ws.Range("C14").Value = 6
$ws.Range("C15").Value = 7

# --- Add new "spi" / TFT LCD section (rows 17-21) ---
$ws.Range("A17").Value = "spi"
$ws.Range("B17").Value = "MISO(RX)"
$ws.Range("C17").Value = 16
$ws.Range("B18").Value = "MOSI(TX)"
$ws.Range("C18").Value = 19
$ws.Range("B19").Value = "SCL"
$ws.Range("C19").Value = 18
$ws.Range("B20").Value = "CS"
$ws.Range("C20").Value = 17
$ws.Range("B21").Value = "D/C"
$ws.Range("C21").Value = 15

# --- View changes: zoom in and move selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 175
$ws.Range("U11").Select()

# --- Drawing: drop the two lower pinout pictures, reposition/resize the remaining one ---
if ($ws.Shapes.Count -ge 3) {
    $ws.Shapes.Item(3).Delete()
}
if ($ws.Shapes.Count -ge 2) {
    $ws.Shapes.Item(2).Delete()
}
$pic = $ws.Shapes.Item(1)
$pic.Left = 279.4472440944882
$pic.Top = 12.96732283464567
$pic.Width = 465.86496062992126
$pic.Height = 282.29976377952755
